$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing "model picture" flag (value 1) for rows 18 and 20,
# matching the pattern used by neighboring rows.
$ws.Range("E18").Value = 1
$ws.Range("E20").Value = 1

# Update the last-selected cell on the sheet to reflect where the
# author ended up after the edit.
$ws.Activate()
$ws.Range("D30").Select() | Out-Null
